$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cultures")
$ws.Range("A36").Value = "Troll"
